$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $escaped = $val.Replace('"', '""')
    $c.Value = '="' + $escaped + '"'
    $c.Copy()
    $c.PasteSpecial(-4163) | Out-Null
}
$excel.CutCopyMode = $false

Set-TextValue "D2" '62.673.96'
Set-TextValue "E2" '  +1.39%  '
Set-TextValue "D3" '2.438.47'
Set-TextValue "E3" '  +1.68%  '
Set-TextValue "E4" '  +0.04%  '
Set-TextValue "D5" '566.40'
Set-TextValue "E5" '  +1.04%  '
Set-TextValue "D6" '145.29'
Set-TextValue "E6" '  +2.40%  '
Set-TextValue "E7" '  -0.05%  '
Set-TextValue "E8" '  +0.21%  '
Set-TextValue "E9" '  +2.36%  '
Set-TextValue "E10" '  +0.34%  '
Set-TextValue "D11" '5.29'
Set-TextValue "E11" '  +1.28%  '
Set-TextValue "D12" '0.355'
Set-TextValue "E12" '  +2.24%  '
Set-TextValue "D13" '26.87'
Set-TextValue "E13" '  +5.83%  '
Set-TextValue "E14" '  +5.87%  '
Set-TextValue "E15" '  +2.06%  '
Set-TextValue "D16" '62.543.67'
Set-TextValue "E16" '  +1.30%  '
Set-TextValue "D17" '2.438.67'
Set-TextValue "E17" '  +1.72%  '
Set-TextValue "E18" '  +0.76%  '
Set-TextValue "E19" '  +2.58%  '
Set-TextValue "D20" '324.09'
Set-TextValue "E20" '  +1.12%  '
Set-TextValue "E21" '  +1.42%  '
Set-TextValue "D22" '0.999'
Set-TextValue "E22" '  -0.13%  '
Set-TextValue "E23" '  +7.20%  '
Set-TextValue "D24" '67.30'
Set-TextValue "E24" '  +2.77%  '
Set-TextValue "E25" '  -1.13%  '
Set-TextValue "D26" '579.65'
Set-TextValue "E26" '  +3.17%  '
Set-TextValue "E27" '  +8.98%  '
Set-TextValue "D29" '0.999'
Set-TextValue "E29" '  -1.50%  '
Set-TextValue "D30" '8.41'
Set-TextValue "E30" '  +3.48%  '
Set-TextValue "E31" '  +4.34%  '
Set-TextValue "E32" '  -0.25%  '
Set-TextValue "E33" '  +0.48%  '
Set-TextValue "E34" '  +0.58%  '
Set-TextValue "E35" '  +2.05%  '
Set-TextValue "E36" '  -0.05%  '
Set-TextValue "D37" '0.382'
Set-TextValue "E37" '  +1.10%  '
Set-TextValue "D38" '18.78'
Set-TextValue "E38" '  +1.95%  '
Set-TextValue "D39" '5.38'
Set-TextValue "E39" '  -0.53%  '
Set-TextValue "D40" '147.66'
Set-TextValue "E40" '  -2.80%  '
Set-TextValue "E41" '  +2.24%  '
Set-TextValue "E42" '  +0.17%  '
Set-TextValue "D43" '2.44'
Set-TextValue "E43" '  +9.66%  '
Set-TextValue "D44" '148.09'
Set-TextValue "E44" '  +0.65%  '
Set-TextValue "E45" '  +2.28%  '
Set-TextValue "D46" '0.0535'
Set-TextValue "E46" '  +1.47%  '
Set-TextValue "D47" '20.52'
Set-TextValue "E47" '  +3.93%  '
Set-TextValue "D48" '0.601'
Set-TextValue "E48" '  +2.67%  '
Set-TextValue "D49" '0.0232'
Set-TextValue "E49" '  +3.45%  '
Set-TextValue "D50" '0.0920'
Set-TextValue "E50" '  +0.79%  '
Set-TextValue "E51" '  +4.92%  '
